$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 46 (columns A:C only, to avoid touching/expanding
# the sheet's used range), shifting existing rows 46+ down by one.
$ws.Range("A46:C46").Insert(-4121)  # xlShiftDown

$ws.Range("A46").Value = "NORMALIZE_DEMAND_TO_ONE"
$ws.Range("B46").Value = $true
$ws.Range("C46").Value = "Normalize demand to 1."

$ws.Range("C46").Select()
